# Instructions_left_hand.xlsx — "Add files via upload" re-upload edit
#
# The authoritative diff shows a handful of worksheet-level data/format
# changes plus a bunch of Excel-version bookkeeping noise (fileVersion
# rupBuild, xr:revisionPtr uids, sharepoint absPath, theme display name,
# unused cellXfs record) that Excel itself rewrites whenever a workbook is
# re-saved by a different Office build and that isn't reachable through the
# documented Application/Workbook/Worksheet/Range object model. This script
# focuses on the concrete, user-visible cell/format edits:
#
#   1. Swap the values in E2/F2 (0.5/0.6 -> 0.6/0.5).
#   2. Drop the (no-op) custom number-format style from E4/F4 so they fall
#      back to the default "Normal" cell style.
#   3. Move the active selection from A4 to F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Swap E2 / F2 values.
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = 0.5

# 2. E4/F4 currently carry a style index that only applies a redundant
#    "General" number format. Resetting to the built-in "Normal" style
#    clears that xf reference so the cells serialize without an s="" index,
#    matching the target (the values themselves are unchanged).
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Style = "Normal"

# 3. Move the selection to F4 (was A4).
$ws.Range("F4").Select()
